$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6 table: switch the table's style (Table Design gallery pick) to
#    a different built-in style id.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tableShape = $s6.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{B2097279-17DA-4D01-B4B6-FE867206D639}")

# ---------------------------------------------------------------------------
# 2) Presentation theme: switch the deck's theme colors from the custom
#    "Integral" palette over to the standard Office palette.
# ---------------------------------------------------------------------------
function Set-ThemeRGB($colors, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $colors.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$tcs = $p.Slides.Item(1).ThemeColorScheme

Set-ThemeRGB $tcs 1  "000000"
Set-ThemeRGB $tcs 2  "FFFFFF"
Set-ThemeRGB $tcs 3  "44546A"
Set-ThemeRGB $tcs 4  "E7E6E6"
Set-ThemeRGB $tcs 5  "5B9BD5"
Set-ThemeRGB $tcs 6  "ED7D31"
Set-ThemeRGB $tcs 7  "A5A5A5"
Set-ThemeRGB $tcs 8  "FFC000"
Set-ThemeRGB $tcs 9  "4472C4"
Set-ThemeRGB $tcs 10 "70AD47"
Set-ThemeRGB $tcs 11 "0563C1"
Set-ThemeRGB $tcs 12 "954F72"
